$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Step 1: remove the two obsolete case rows (A 42417-2025 and A 42412-2025), originally rows 29-30
$ws.Rows(29).Delete() | Out-Null
$ws.Rows(29).Delete() | Out-Null

# Step 2: clear the old F15 value ("Kommuner") - it moves to a different row after the reorder
$ws.Range("F15").ClearContents() | Out-Null

# Step 3: rewrite rows 2-33 with the reordered/updated case data

# Row 2: A 2943-2023
$ws.Range("A2").Value = "A 2943-2023"
$ws.Range("B2").Value = 44945
$ws.Range("C2").Value = 46065
$ws.Range("D2").Value = "SKÅNE LÄN"
$ws.Range("E2").Value = "BÅSTAD"
$ws.Range("G2").Value = 11.3
$ws.Range("H2").Value = 1
$ws.Range("I2").Value = 1
$ws.Range("J2").Value = 0
$ws.Range("K2").Value = 0
$ws.Range("L2").Value = 0
$ws.Range("M2").Value = 1
$ws.Range("N2").Value = 0
$ws.Range("O2").Value = 1
$ws.Range("P2").Value = 1
$ws.Range("Q2").Value = 3
$ws.Range("R2").Value = "Skogsalm`r`nSårläka`r`nGrönvit nattviol"
$ws.Range("S2").Formula = "=HYPERLINK(""https://klasma.github.io/Logging_1278/artfynd/A 2943-2023 artfynd.xlsx"", ""A 2943-2023"")"
$ws.Range("T2").Formula = "=HYPERLINK(""https://klasma.github.io/Logging_1278/kartor/A 2943-2023 karta.png"", ""A 2943-2023"")"
$ws.Range("V2").Formula = "=HYPERLINK(""https://klasma.github.io/Logging_1278/klagomål/A 2943-2023 FSC-klagomål.docx"", ""A 2943-2023"")"
$ws.Range("W2").Formula = "=HYPERLINK(""https://klasma.github.io/Logging_1278/klagomålsmail/A 2943-2023 FSC-klagomål mail.docx"", ""A 2943-2023"")"
$ws.Range("X2").Formula = "=HYPERLINK(""https://klasma.github.io/Logging_1278/tillsyn/A 2943-2023 tillsynsbegäran.docx"", ""A 2943-2023"")"
$ws.Range("Y2").Formula = "=HYPERLINK(""https://klasma.github.io/Logging_1278/tillsynsmail/A 2943-2023 tillsynsbegäran mail.docx"", ""A 2943-2023"")"

# Row 3: A 18319-2025
$ws.Range("A3").Value = "A 18319-2025"
$ws.Range("B3").Value = 45762
$ws.Range("C3").Value = 46065
$ws.Range("D3").Value = "SKÅNE LÄN"
$ws.Range("E3").Value = "BÅSTAD"
$ws.Range("G3").Value = 3.4
$ws.Range("H3").Value = 1
$ws.Range("I3").Value = 2
$ws.Range("J3").Value = 1
$ws.Range("K3").Value = 0
$ws.Range("L3").Value = 0
$ws.Range("M3").Value = 0
$ws.Range("N3").Value = 0
$ws.Range("O3").Value = 1
$ws.Range("P3").Value = 0
$ws.Range("Q3").Value = 3
$ws.Range("R3").Value = "Igelkott`r`nKorallrot`r`nKällpraktmossa"
$ws.Range("S3").Formula = "=HYPERLINK(""https://klasma.github.io/Logging_1278/artfynd/A 18319-2025 artfynd.xlsx"", ""A 18319-2025"")"
$ws.Range("T3").Formula = "=HYPERLINK(""https://klasma.github.io/Logging_1278/kartor/A 18319-2025 karta.png"", ""A 18319-2025"")"
$ws.Range("V3").Formula = "=HYPERLINK(""https://klasma.github.io/Logging_1278/klagomål/A 18319-2025 FSC-klagomål.docx"", ""A 18319-2025"")"
$ws.Range("W3").Formula = "=HYPERLINK(""https://klasma.github.io/Logging_1278/klagomålsmail/A 18319-2025 FSC-klagomål mail.docx"", ""A 18319-2025"")"
$ws.Range("X3").Formula = "=HYPERLINK(""https://klasma.github.io/Logging_1278/tillsyn/A 18319-2025 tillsynsbegäran.docx"", ""A 18319-2025"")"
$ws.Range("Y3").Formula = "=HYPERLINK(""https://klasma.github.io/Logging_1278/tillsynsmail/A 18319-2025 tillsynsbegäran mail.docx"", ""A 18319-2025"")"

# Row 4: A 1782-2024
$ws.Range("A4").Value = "A 1782-2024"
$ws.Range("B4").Value = 45307
$ws.Range("C4").Value = 46065
$ws.Range("D4").Value = "SKÅNE LÄN"
$ws.Range("E4").Value = "BÅSTAD"
$ws.Range("G4").Value = 2.7
$ws.Range("H4").Value = 1
$ws.Range("I4").Value = 0
$ws.Range("J4").Value = 0
$ws.Range("K4").Value = 0
$ws.Range("L4").Value = 0
$ws.Range("M4").Value = 0
$ws.Range("N4").Value = 0
$ws.Range("O4").Value = 0
$ws.Range("P4").Value = 0
$ws.Range("Q4").Value = 1
$ws.Range("R4").Value = "Grönvit nattviol"
$ws.Range("S4").Formula = "=HYPERLINK(""https://klasma.github.io/Logging_1278/artfynd/A 1782-2024 artfynd.xlsx"", ""A 1782-2024"")"
$ws.Range("T4").Formula = "=HYPERLINK(""https://klasma.github.io/Logging_1278/kartor/A 1782-2024 karta.png"", ""A 1782-2024"")"
$ws.Range("V4").Formula = "=HYPERLINK(""https://klasma.github.io/Logging_1278/klagomål/A 1782-2024 FSC-klagomål.docx"", ""A 1782-2024"")"
$ws.Range("W4").Formula = "=HYPERLINK(""https://klasma.github.io/Logging_1278/klagomålsmail/A 1782-2024 FSC-klagomål mail.docx"", ""A 1782-2024"")"
$ws.Range("X4").Formula = "=HYPERLINK(""https://klasma.github.io/Logging_1278/tillsyn/A 1782-2024 tillsynsbegäran.docx"", ""A 1782-2024"")"
$ws.Range("Y4").Formula = "=HYPERLINK(""https://klasma.github.io/Logging_1278/tillsynsmail/A 1782-2024 tillsynsbegäran mail.docx"", ""A 1782-2024"")"

# Row 5: A 389-2023
$ws.Range("A5").Value = "A 389-2023"
$ws.Range("B5").Value = 44929
$ws.Range("C5").Value = 46065
$ws.Range("D5").Value = "SKÅNE LÄN"
$ws.Range("E5").Value = "BÅSTAD"
$ws.Range("G5").Value = 2.5
$ws.Range("H5").Value = 1
$ws.Range("I5").Value = 0
$ws.Range("J5").Value = 0
$ws.Range("K5").Value = 0
$ws.Range("L5").Value = 0
$ws.Range("M5").Value = 0
$ws.Range("N5").Value = 0
$ws.Range("O5").Value = 0
$ws.Range("P5").Value = 0
$ws.Range("Q5").Value = 1
$ws.Range("R5").Value = "Grönvit nattviol"
$ws.Range("S5").Formula = "=HYPERLINK(""https://klasma.github.io/Logging_1278/artfynd/A 389-2023 artfynd.xlsx"", ""A 389-2023"")"
$ws.Range("T5").Formula = "=HYPERLINK(""https://klasma.github.io/Logging_1278/kartor/A 389-2023 karta.png"", ""A 389-2023"")"
$ws.Range("V5").Formula = "=HYPERLINK(""https://klasma.github.io/Logging_1278/klagomål/A 389-2023 FSC-klagomål.docx"", ""A 389-2023"")"
$ws.Range("W5").Formula = "=HYPERLINK(""https://klasma.github.io/Logging_1278/klagomålsmail/A 389-2023 FSC-klagomål mail.docx"", ""A 389-2023"")"
$ws.Range("X5").Formula = "=HYPERLINK(""https://klasma.github.io/Logging_1278/tillsyn/A 389-2023 tillsynsbegäran.docx"", ""A 389-2023"")"
$ws.Range("Y5").Formula = "=HYPERLINK(""https://klasma.github.io/Logging_1278/tillsynsmail/A 389-2023 tillsynsbegäran mail.docx"", ""A 389-2023"")"

# Row 6: A 38013-2022
$ws.Range("A6").Value = "A 38013-2022"
$ws.Range("B6").Value = 44811
$ws.Range("C6").Value = 46065
$ws.Range("D6").Value = "SKÅNE LÄN"
$ws.Range("E6").Value = "BÅSTAD"
$ws.Range("G6").Value = 5.3
$ws.Range("H6").Value = 0
$ws.Range("I6").Value = 0
$ws.Range("J6").Value = 1
$ws.Range("K6").Value = 0
$ws.Range("L6").Value = 0
$ws.Range("M6").Value = 0
$ws.Range("N6").Value = 0
$ws.Range("O6").Value = 1
$ws.Range("P6").Value = 0
$ws.Range("Q6").Value = 1
$ws.Range("R6").Value = "Igelkott"
$ws.Range("S6").Formula = "=HYPERLINK(""https://klasma.github.io/Logging_1278/artfynd/A 38013-2022 artfynd.xlsx"", ""A 38013-2022"")"
$ws.Range("T6").Formula = "=HYPERLINK(""https://klasma.github.io/Logging_1278/kartor/A 38013-2022 karta.png"", ""A 38013-2022"")"
$ws.Range("V6").Formula = "=HYPERLINK(""https://klasma.github.io/Logging_1278/klagomål/A 38013-2022 FSC-klagomål.docx"", ""A 38013-2022"")"
$ws.Range("W6").Formula = "=HYPERLINK(""https://klasma.github.io/Logging_1278/klagomålsmail/A 38013-2022 FSC-klagomål mail.docx"", ""A 38013-2022"")"
$ws.Range("X6").Formula = "=HYPERLINK(""https://klasma.github.io/Logging_1278/tillsyn/A 38013-2022 tillsynsbegäran.docx"", ""A 38013-2022"")"
$ws.Range("Y6").Formula = "=HYPERLINK(""https://klasma.github.io/Logging_1278/tillsynsmail/A 38013-2022 tillsynsbegäran mail.docx"", ""A 38013-2022"")"

# Row 7: A 29245-2021
$ws.Range("A7").Value = "A 29245-2021"
$ws.Range("B7").Value = 44361
$ws.Range("C7").Value = 46065
$ws.Range("D7").Value = "SKÅNE LÄN"
$ws.Range("E7").Value = "BÅSTAD"
$ws.Range("G7").Value = 0.4
$ws.Range("H7").Value = 0
$ws.Range("I7").Value = 0
$ws.Range("J7").Value = 0
$ws.Range("K7").Value = 0
$ws.Range("L7").Value = 0
$ws.Range("M7").Value = 0
$ws.Range("N7").Value = 0
$ws.Range("O7").Value = 0
$ws.Range("P7").Value = 0
$ws.Range("Q7").Value = 0

# Row 8: A 24-2023
$ws.Range("A8").Value = "A 24-2023"
$ws.Range("B8").Value = 44928
$ws.Range("C8").Value = 46065
$ws.Range("D8").Value = "SKÅNE LÄN"
$ws.Range("E8").Value = "BÅSTAD"
$ws.Range("G8").Value = 0.5
$ws.Range("H8").Value = 0
$ws.Range("I8").Value = 0
$ws.Range("J8").Value = 0
$ws.Range("K8").Value = 0
$ws.Range("L8").Value = 0
$ws.Range("M8").Value = 0
$ws.Range("N8").Value = 0
$ws.Range("O8").Value = 0
$ws.Range("P8").Value = 0
$ws.Range("Q8").Value = 0

# Row 9: A 21572-2023
$ws.Range("A9").Value = "A 21572-2023"
$ws.Range("B9").Value = 45063
$ws.Range("C9").Value = 46065
$ws.Range("D9").Value = "SKÅNE LÄN"
$ws.Range("E9").Value = "BÅSTAD"
$ws.Range("G9").Value = 1.7
$ws.Range("H9").Value = 0
$ws.Range("I9").Value = 0
$ws.Range("J9").Value = 0
$ws.Range("K9").Value = 0
$ws.Range("L9").Value = 0
$ws.Range("M9").Value = 0
$ws.Range("N9").Value = 0
$ws.Range("O9").Value = 0
$ws.Range("P9").Value = 0
$ws.Range("Q9").Value = 0

# Row 10: A 1531-2022
$ws.Range("A10").Value = "A 1531-2022"
$ws.Range("B10").Value = 44573
$ws.Range("C10").Value = 46065
$ws.Range("D10").Value = "SKÅNE LÄN"
$ws.Range("E10").Value = "BÅSTAD"
$ws.Range("G10").Value = 1.6
$ws.Range("H10").Value = 0
$ws.Range("I10").Value = 0
$ws.Range("J10").Value = 0
$ws.Range("K10").Value = 0
$ws.Range("L10").Value = 0
$ws.Range("M10").Value = 0
$ws.Range("N10").Value = 0
$ws.Range("O10").Value = 0
$ws.Range("P10").Value = 0
$ws.Range("Q10").Value = 0

# Row 11: A 4486-2024
$ws.Range("A11").Value = "A 4486-2024"
$ws.Range("B11").Value = 45327
$ws.Range("C11").Value = 46065
$ws.Range("D11").Value = "SKÅNE LÄN"
$ws.Range("E11").Value = "BÅSTAD"
$ws.Range("G11").Value = 0.6
$ws.Range("H11").Value = 0
$ws.Range("I11").Value = 0
$ws.Range("J11").Value = 0
$ws.Range("K11").Value = 0
$ws.Range("L11").Value = 0
$ws.Range("M11").Value = 0
$ws.Range("N11").Value = 0
$ws.Range("O11").Value = 0
$ws.Range("P11").Value = 0
$ws.Range("Q11").Value = 0

# Row 12: A 48974-2023
$ws.Range("A12").Value = "A 48974-2023"
$ws.Range("B12").Value = 45209
$ws.Range("C12").Value = 46065
$ws.Range("D12").Value = "SKÅNE LÄN"
$ws.Range("E12").Value = "BÅSTAD"
$ws.Range("G12").Value = 4.5
$ws.Range("H12").Value = 0
$ws.Range("I12").Value = 0
$ws.Range("J12").Value = 0
$ws.Range("K12").Value = 0
$ws.Range("L12").Value = 0
$ws.Range("M12").Value = 0
$ws.Range("N12").Value = 0
$ws.Range("O12").Value = 0
$ws.Range("P12").Value = 0
$ws.Range("Q12").Value = 0

# Row 13: A 4822-2023
$ws.Range("A13").Value = "A 4822-2023"
$ws.Range("B13").Value = 44957
$ws.Range("C13").Value = 46065
$ws.Range("D13").Value = "SKÅNE LÄN"
$ws.Range("E13").Value = "BÅSTAD"
$ws.Range("G13").Value = 2.2
$ws.Range("H13").Value = 0
$ws.Range("I13").Value = 0
$ws.Range("J13").Value = 0
$ws.Range("K13").Value = 0
$ws.Range("L13").Value = 0
$ws.Range("M13").Value = 0
$ws.Range("N13").Value = 0
$ws.Range("O13").Value = 0
$ws.Range("P13").Value = 0
$ws.Range("Q13").Value = 0

# Row 14: A 32610-2024
$ws.Range("A14").Value = "A 32610-2024"
$ws.Range("B14").Value = 45513
$ws.Range("C14").Value = 46065
$ws.Range("D14").Value = "SKÅNE LÄN"
$ws.Range("E14").Value = "BÅSTAD"
$ws.Range("G14").Value = 0.5
$ws.Range("H14").Value = 0
$ws.Range("I14").Value = 0
$ws.Range("J14").Value = 0
$ws.Range("K14").Value = 0
$ws.Range("L14").Value = 0
$ws.Range("M14").Value = 0
$ws.Range("N14").Value = 0
$ws.Range("O14").Value = 0
$ws.Range("P14").Value = 0
$ws.Range("Q14").Value = 0

# Row 15: A 28260-2023
$ws.Range("A15").Value = "A 28260-2023"
$ws.Range("B15").Value = 45099
$ws.Range("C15").Value = 46065
$ws.Range("D15").Value = "SKÅNE LÄN"
$ws.Range("E15").Value = "BÅSTAD"
$ws.Range("G15").Value = 5
$ws.Range("H15").Value = 0
$ws.Range("I15").Value = 0
$ws.Range("J15").Value = 0
$ws.Range("K15").Value = 0
$ws.Range("L15").Value = 0
$ws.Range("M15").Value = 0
$ws.Range("N15").Value = 0
$ws.Range("O15").Value = 0
$ws.Range("P15").Value = 0
$ws.Range("Q15").Value = 0

# Row 16: A 635-2023
$ws.Range("A16").Value = "A 635-2023"
$ws.Range("B16").Value = 44930
$ws.Range("C16").Value = 46065
$ws.Range("D16").Value = "SKÅNE LÄN"
$ws.Range("E16").Value = "BÅSTAD"
$ws.Range("G16").Value = 0.5
$ws.Range("H16").Value = 0
$ws.Range("I16").Value = 0
$ws.Range("J16").Value = 0
$ws.Range("K16").Value = 0
$ws.Range("L16").Value = 0
$ws.Range("M16").Value = 0
$ws.Range("N16").Value = 0
$ws.Range("O16").Value = 0
$ws.Range("P16").Value = 0
$ws.Range("Q16").Value = 0

# Row 17: A 48181-2024
$ws.Range("A17").Value = "A 48181-2024"
$ws.Range("B17").Value = 45589
$ws.Range("C17").Value = 46065
$ws.Range("D17").Value = "SKÅNE LÄN"
$ws.Range("E17").Value = "BÅSTAD"
$ws.Range("G17").Value = 0.7
$ws.Range("H17").Value = 0
$ws.Range("I17").Value = 0
$ws.Range("J17").Value = 0
$ws.Range("K17").Value = 0
$ws.Range("L17").Value = 0
$ws.Range("M17").Value = 0
$ws.Range("N17").Value = 0
$ws.Range("O17").Value = 0
$ws.Range("P17").Value = 0
$ws.Range("Q17").Value = 0

# Row 18: A 4481-2024
$ws.Range("A18").Value = "A 4481-2024"
$ws.Range("B18").Value = 45327
$ws.Range("C18").Value = 46065
$ws.Range("D18").Value = "SKÅNE LÄN"
$ws.Range("E18").Value = "BÅSTAD"
$ws.Range("G18").Value = 1
$ws.Range("H18").Value = 0
$ws.Range("I18").Value = 0
$ws.Range("J18").Value = 0
$ws.Range("K18").Value = 0
$ws.Range("L18").Value = 0
$ws.Range("M18").Value = 0
$ws.Range("N18").Value = 0
$ws.Range("O18").Value = 0
$ws.Range("P18").Value = 0
$ws.Range("Q18").Value = 0

# Row 19: A 10710-2025
$ws.Range("A19").Value = "A 10710-2025"
$ws.Range("B19").Value = 45722
$ws.Range("C19").Value = 46065
$ws.Range("D19").Value = "SKÅNE LÄN"
$ws.Range("E19").Value = "BÅSTAD"
$ws.Range("F19").Value = "Kommuner"
$ws.Range("G19").Value = 1.8
$ws.Range("H19").Value = 0
$ws.Range("I19").Value = 0
$ws.Range("J19").Value = 0
$ws.Range("K19").Value = 0
$ws.Range("L19").Value = 0
$ws.Range("M19").Value = 0
$ws.Range("N19").Value = 0
$ws.Range("O19").Value = 0
$ws.Range("P19").Value = 0
$ws.Range("Q19").Value = 0

# Row 20: A 18328-2025
$ws.Range("A20").Value = "A 18328-2025"
$ws.Range("B20").Value = 45762
$ws.Range("C20").Value = 46065
$ws.Range("D20").Value = "SKÅNE LÄN"
$ws.Range("E20").Value = "BÅSTAD"
$ws.Range("G20").Value = 1.8
$ws.Range("H20").Value = 0
$ws.Range("I20").Value = 0
$ws.Range("J20").Value = 0
$ws.Range("K20").Value = 0
$ws.Range("L20").Value = 0
$ws.Range("M20").Value = 0
$ws.Range("N20").Value = 0
$ws.Range("O20").Value = 0
$ws.Range("P20").Value = 0
$ws.Range("Q20").Value = 0

# Row 21: A 18332-2025
$ws.Range("A21").Value = "A 18332-2025"
$ws.Range("B21").Value = 45762
$ws.Range("C21").Value = 46065
$ws.Range("D21").Value = "SKÅNE LÄN"
$ws.Range("E21").Value = "BÅSTAD"
$ws.Range("G21").Value = 2.5
$ws.Range("H21").Value = 0
$ws.Range("I21").Value = 0
$ws.Range("J21").Value = 0
$ws.Range("K21").Value = 0
$ws.Range("L21").Value = 0
$ws.Range("M21").Value = 0
$ws.Range("N21").Value = 0
$ws.Range("O21").Value = 0
$ws.Range("P21").Value = 0
$ws.Range("Q21").Value = 0

# Row 22: A 34400-2025
$ws.Range("A22").Value = "A 34400-2025"
$ws.Range("B22").Value = 45846.61351851852
$ws.Range("C22").Value = 46065
$ws.Range("D22").Value = "SKÅNE LÄN"
$ws.Range("E22").Value = "BÅSTAD"
$ws.Range("G22").Value = 1.3
$ws.Range("H22").Value = 0
$ws.Range("I22").Value = 0
$ws.Range("J22").Value = 0
$ws.Range("K22").Value = 0
$ws.Range("L22").Value = 0
$ws.Range("M22").Value = 0
$ws.Range("N22").Value = 0
$ws.Range("O22").Value = 0
$ws.Range("P22").Value = 0
$ws.Range("Q22").Value = 0

# Row 23: A 34401-2025
$ws.Range("A23").Value = "A 34401-2025"
$ws.Range("B23").Value = 45846.6140162037
$ws.Range("C23").Value = 46065
$ws.Range("D23").Value = "SKÅNE LÄN"
$ws.Range("E23").Value = "BÅSTAD"
$ws.Range("G23").Value = 2.8
$ws.Range("H23").Value = 0
$ws.Range("I23").Value = 0
$ws.Range("J23").Value = 0
$ws.Range("K23").Value = 0
$ws.Range("L23").Value = 0
$ws.Range("M23").Value = 0
$ws.Range("N23").Value = 0
$ws.Range("O23").Value = 0
$ws.Range("P23").Value = 0
$ws.Range("Q23").Value = 0

# Row 24: A 4487-2024
$ws.Range("A24").Value = "A 4487-2024"
$ws.Range("B24").Value = 45327
$ws.Range("C24").Value = 46065
$ws.Range("D24").Value = "SKÅNE LÄN"
$ws.Range("E24").Value = "BÅSTAD"
$ws.Range("G24").Value = 1.9
$ws.Range("H24").Value = 0
$ws.Range("I24").Value = 0
$ws.Range("J24").Value = 0
$ws.Range("K24").Value = 0
$ws.Range("L24").Value = 0
$ws.Range("M24").Value = 0
$ws.Range("N24").Value = 0
$ws.Range("O24").Value = 0
$ws.Range("P24").Value = 0
$ws.Range("Q24").Value = 0

# Row 25: A 4256-2025
$ws.Range("A25").Value = "A 4256-2025"
$ws.Range("B25").Value = 45685
$ws.Range("C25").Value = 46065
$ws.Range("D25").Value = "SKÅNE LÄN"
$ws.Range("E25").Value = "BÅSTAD"
$ws.Range("G25").Value = 2
$ws.Range("H25").Value = 0
$ws.Range("I25").Value = 0
$ws.Range("J25").Value = 0
$ws.Range("K25").Value = 0
$ws.Range("L25").Value = 0
$ws.Range("M25").Value = 0
$ws.Range("N25").Value = 0
$ws.Range("O25").Value = 0
$ws.Range("P25").Value = 0
$ws.Range("Q25").Value = 0

# Row 26: A 5817-2025
$ws.Range("A26").Value = "A 5817-2025"
$ws.Range("B26").Value = 45694.74113425926
$ws.Range("C26").Value = 46065
$ws.Range("D26").Value = "SKÅNE LÄN"
$ws.Range("E26").Value = "BÅSTAD"
$ws.Range("G26").Value = 1.2
$ws.Range("H26").Value = 0
$ws.Range("I26").Value = 0
$ws.Range("J26").Value = 0
$ws.Range("K26").Value = 0
$ws.Range("L26").Value = 0
$ws.Range("M26").Value = 0
$ws.Range("N26").Value = 0
$ws.Range("O26").Value = 0
$ws.Range("P26").Value = 0
$ws.Range("Q26").Value = 0

# Row 27: A 11517-2024
$ws.Range("A27").Value = "A 11517-2024"
$ws.Range("B27").Value = 45372
$ws.Range("C27").Value = 46065
$ws.Range("D27").Value = "SKÅNE LÄN"
$ws.Range("E27").Value = "BÅSTAD"
$ws.Range("G27").Value = 0.7
$ws.Range("H27").Value = 0
$ws.Range("I27").Value = 0
$ws.Range("J27").Value = 0
$ws.Range("K27").Value = 0
$ws.Range("L27").Value = 0
$ws.Range("M27").Value = 0
$ws.Range("N27").Value = 0
$ws.Range("O27").Value = 0
$ws.Range("P27").Value = 0
$ws.Range("Q27").Value = 0

# Row 28: A 18434-2023
$ws.Range("A28").Value = "A 18434-2023"
$ws.Range("B28").Value = 45042
$ws.Range("C28").Value = 46065
$ws.Range("D28").Value = "SKÅNE LÄN"
$ws.Range("E28").Value = "BÅSTAD"
$ws.Range("G28").Value = 0.7
$ws.Range("H28").Value = 0
$ws.Range("I28").Value = 0
$ws.Range("J28").Value = 0
$ws.Range("K28").Value = 0
$ws.Range("L28").Value = 0
$ws.Range("M28").Value = 0
$ws.Range("N28").Value = 0
$ws.Range("O28").Value = 0
$ws.Range("P28").Value = 0
$ws.Range("Q28").Value = 0

# Row 29: A 53131-2021
$ws.Range("A29").Value = "A 53131-2021"
$ws.Range("B29").Value = 44468
$ws.Range("C29").Value = 46065
$ws.Range("D29").Value = "SKÅNE LÄN"
$ws.Range("E29").Value = "BÅSTAD"
$ws.Range("G29").Value = 1.3
$ws.Range("H29").Value = 0
$ws.Range("I29").Value = 0
$ws.Range("J29").Value = 0
$ws.Range("K29").Value = 0
$ws.Range("L29").Value = 0
$ws.Range("M29").Value = 0
$ws.Range("N29").Value = 0
$ws.Range("O29").Value = 0
$ws.Range("P29").Value = 0
$ws.Range("Q29").Value = 0

# Row 30: A 18327-2025
$ws.Range("A30").Value = "A 18327-2025"
$ws.Range("B30").Value = 45762
$ws.Range("C30").Value = 46065
$ws.Range("D30").Value = "SKÅNE LÄN"
$ws.Range("E30").Value = "BÅSTAD"
$ws.Range("G30").Value = 0.6
$ws.Range("H30").Value = 0
$ws.Range("I30").Value = 0
$ws.Range("J30").Value = 0
$ws.Range("K30").Value = 0
$ws.Range("L30").Value = 0
$ws.Range("M30").Value = 0
$ws.Range("N30").Value = 0
$ws.Range("O30").Value = 0
$ws.Range("P30").Value = 0
$ws.Range("Q30").Value = 0

# Row 31: A 4493-2024
$ws.Range("A31").Value = "A 4493-2024"
$ws.Range("B31").Value = 45327
$ws.Range("C31").Value = 46065
$ws.Range("D31").Value = "SKÅNE LÄN"
$ws.Range("E31").Value = "BÅSTAD"
$ws.Range("G31").Value = 1.8
$ws.Range("H31").Value = 0
$ws.Range("I31").Value = 0
$ws.Range("J31").Value = 0
$ws.Range("K31").Value = 0
$ws.Range("L31").Value = 0
$ws.Range("M31").Value = 0
$ws.Range("N31").Value = 0
$ws.Range("O31").Value = 0
$ws.Range("P31").Value = 0
$ws.Range("Q31").Value = 0

# Row 32: A 7731-2026
$ws.Range("A32").Value = "A 7731-2026"
$ws.Range("B32").Value = 46062.52008101852
$ws.Range("C32").Value = 46065
$ws.Range("D32").Value = "SKÅNE LÄN"
$ws.Range("E32").Value = "BÅSTAD"
$ws.Range("G32").Value = 5.9
$ws.Range("H32").Value = 0
$ws.Range("I32").Value = 0
$ws.Range("J32").Value = 0
$ws.Range("K32").Value = 0
$ws.Range("L32").Value = 0
$ws.Range("M32").Value = 0
$ws.Range("N32").Value = 0
$ws.Range("O32").Value = 0
$ws.Range("P32").Value = 0
$ws.Range("Q32").Value = 0

# Row 33: A 7727-2026
$ws.Range("A33").Value = "A 7727-2026"
$ws.Range("B33").Value = 46062.50420138889
$ws.Range("C33").Value = 46065
$ws.Range("D33").Value = "SKÅNE LÄN"
$ws.Range("E33").Value = "BÅSTAD"
$ws.Range("G33").Value = 1.9
$ws.Range("H33").Value = 0
$ws.Range("I33").Value = 0
$ws.Range("J33").Value = 0
$ws.Range("K33").Value = 0
$ws.Range("L33").Value = 0
$ws.Range("M33").Value = 0
$ws.Range("N33").Value = 0
$ws.Range("O33").Value = 0
$ws.Range("P33").Value = 0
$ws.Range("Q33").Value = 0
